$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value2 = 1969.2307  # H40
$ws.Cells.Item(40, 9).Value2 = 1936.8422  # I40
$ws.Cells.Item(40, 11).Value2 = 1936.8422  # K40
$ws.Cells.Item(40, 13).Value2 = -1761.8422  # M40
$ws.Cells.Item(69, 8).Value2 = 4307.5  # H69
$ws.Cells.Item(69, 9).Value2 = 4000  # I69
$ws.Cells.Item(69, 11).Value2 = 12000  # K69
$ws.Cells.Item(69, 13).Value2 = -11126  # M69
$ws.Cells.Item(72, 8).Value2 = 4307.5  # H72
$ws.Cells.Item(72, 9).Value2 = 4000  # I72
$ws.Cells.Item(72, 11).Value2 = 36000  # K72
$ws.Cells.Item(72, 13).Value2 = -31632  # M72
$ws.Cells.Item(132, 8).Value2 = 7696763  # H132
$ws.Cells.Item(132, 9).Value2 = 10003742  # I132
$ws.Cells.Item(132, 10).Value2 = 6833.3335  # J132
$ws.Cells.Item(132, 11).Value2 = 30011226  # K132
$ws.Cells.Item(132, 12).Value2 = 20500.0005  # L132
$ws.Cells.Item(132, 13).Value2 = -30008696  # M132
$ws.Cells.Item(132, 14).Value2 = -25560.0005  # N132
$ws.Cells.Item(137, 8).Value2 = 3594.3555  # H137
$ws.Cells.Item(137, 9).Value2 = 4064.8965  # I137
$ws.Cells.Item(137, 10).Value2 = 2741.5  # J137
$ws.Cells.Item(137, 11).Value2 = 12194.6895  # K137
$ws.Cells.Item(137, 12).Value2 = 8224.5  # L137
$ws.Cells.Item(137, 13).Value2 = -9644.6895  # M137
$ws.Cells.Item(137, 14).Value2 = -13324.5  # N137
$ws.Cells.Item(138, 8).Value2 = 4613.28  # H138
$ws.Cells.Item(138, 9).Value2 = 2337.2307  # I138
$ws.Cells.Item(138, 10).Value2 = 5820.9795  # J138
$ws.Cells.Item(138, 11).Value2 = 7011.6921  # K138
$ws.Cells.Item(138, 12).Value2 = 17462.9385  # L138
$ws.Cells.Item(138, 13).Value2 = -1871.6921  # M138
$ws.Cells.Item(138, 14).Value2 = -27742.9385  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value2 = 74105.39999999999  # H16
$ws.Cells.Item(16, 9).Value2 = 500  # I16
$ws.Cells.Item(16, 10).Value2 = 92506.75  # J16
$ws.Cells.Item(16, 11).Value2 = 500  # K16
$ws.Cells.Item(16, 12).Value2 = 92506.75  # L16
$ws.Cells.Item(16, 13).Value2 = -213  # M16
$ws.Cells.Item(16, 14).Value2 = -93080.75  # N16
$ws.Cells.Item(32, 8).Value2 = 2884.0918  # H32
$ws.Cells.Item(32, 9).Value2 = 2884.0918  # I32
$ws.Cells.Item(32, 11).Value2 = 2884.0918  # K32
$ws.Cells.Item(32, 13).Value2 = -2597.0918  # M32
$ws.Cells.Item(60, 8).Value2 = 10620  # H60
$ws.Cells.Item(60, 9).Value2 = 6183  # I60
$ws.Cells.Item(60, 10).Value2 = 15057  # J60
$ws.Cells.Item(60, 11).Value2 = 6183  # K60
$ws.Cells.Item(60, 12).Value2 = 15057  # L60
$ws.Cells.Item(60, 13).Value2 = -5450  # M60
$ws.Cells.Item(60, 14).Value2 = -16523  # N60
$ws.Cells.Item(117, 8).Value2 = 29666.666  # H117
$ws.Cells.Item(117, 10).Value2 = 29666.666  # J117
$ws.Cells.Item(117, 12).Value2 = 29666.666  # L117
$ws.Cells.Item(117, 14).Value2 = -38844.666  # N117

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(122, 8).Value2 = 30000  # H122
$ws.Cells.Item(122, 10).Value2 = 30000  # J122
$ws.Cells.Item(122, 12).Value2 = 30000  # L122
$ws.Cells.Item(122, 14).Value2 = -39800  # N122
$ws.Cells.Item(134, 8).Value2 = 3457.96  # H134
$ws.Cells.Item(134, 9).Value2 = 3185.375  # I134
$ws.Cells.Item(134, 10).Value2 = 10000  # J134
$ws.Cells.Item(134, 11).Value2 = 9556.125  # K134
$ws.Cells.Item(134, 12).Value2 = 30000  # L134
$ws.Cells.Item(134, 13).Value2 = -7021.125  # M134
$ws.Cells.Item(134, 14).Value2 = -35070  # N134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value2 = 58756.125  # H23
$ws.Cells.Item(23, 10).Value2 = 65720  # J23
$ws.Cells.Item(23, 12).Value2 = 65720  # L23
$ws.Cells.Item(23, 14).Value2 = -66200  # N23
$ws.Cells.Item(27, 8).Value2 = 58756.125  # H27
$ws.Cells.Item(27, 10).Value2 = 65720  # J27
$ws.Cells.Item(27, 12).Value2 = 65720  # L27
$ws.Cells.Item(27, 14).Value2 = -66104  # N27
$ws.Cells.Item(116, 8).Value2 = 32800  # H116
$ws.Cells.Item(116, 10).Value2 = 32800  # J116
$ws.Cells.Item(116, 12).Value2 = 32800  # L116
$ws.Cells.Item(116, 14).Value2 = -41978  # N116
$ws.Cells.Item(132, 8).Value2 = 2271.5  # H132
$ws.Cells.Item(132, 9).Value2 = 1782.5333  # I132
$ws.Cells.Item(132, 10).Value2 = 4716.3335  # J132
$ws.Cells.Item(132, 11).Value2 = 5347.5999  # K132
$ws.Cells.Item(132, 12).Value2 = 14149.0005  # L132
$ws.Cells.Item(132, 13).Value2 = -2817.5999  # M132
$ws.Cells.Item(132, 14).Value2 = -19209.0005  # N132
$ws.Cells.Item(134, 8).Value2 = 10418774  # H134
$ws.Cells.Item(134, 9).Value2 = 13515161  # I134
$ws.Cells.Item(134, 10).Value2 = 3654.4546  # J134
$ws.Cells.Item(134, 11).Value2 = 40545483  # K134
$ws.Cells.Item(134, 12).Value2 = 10963.3638  # L134
$ws.Cells.Item(134, 13).Value2 = -40542948  # M134
$ws.Cells.Item(134, 14).Value2 = -16033.3638  # N134

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value2 = 296.6  # H14
$ws.Cells.Item(14, 9).Value2 = 296.6  # I14
$ws.Cells.Item(14, 11).Value2 = 889.8000000000001  # K14
$ws.Cells.Item(14, 13).Value2 = -716.8000000000001  # M14
$ws.Cells.Item(56, 8).Value2 = 5279.2856  # H56
$ws.Cells.Item(56, 9).Value2 = 5279.2856  # I56
$ws.Cells.Item(56, 11).Value2 = 5279.2856  # K56
$ws.Cells.Item(56, 13).Value2 = -4749.2856  # M56
$ws.Cells.Item(87, 8).Value2 = 14283.333  # H87
$ws.Cells.Item(87, 10).Value2 = 15950  # J87
$ws.Cells.Item(87, 12).Value2 = 47850  # L87
$ws.Cells.Item(87, 14).Value2 = -50346  # N87
$ws.Cells.Item(90, 8).Value2 = 14283.333  # H90
$ws.Cells.Item(90, 10).Value2 = 15950  # J90
$ws.Cells.Item(90, 12).Value2 = 143550  # L90
$ws.Cells.Item(90, 14).Value2 = -156030  # N90

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value2 = 3122267  # H12
$ws.Cells.Item(12, 9).Value2 = 2579010.5  # I12
$ws.Cells.Item(12, 11).Value2 = 2579010.5  # K12
$ws.Cells.Item(12, 13).Value2 = -2578870.5  # M12
$ws.Cells.Item(80, 8).Value2 = 3165.5  # H80
$ws.Cells.Item(80, 9).Value2 = 2784.074  # I80
$ws.Cells.Item(80, 10).Value2 = 3771.2942  # J80
$ws.Cells.Item(80, 11).Value2 = 2784.074  # K80
$ws.Cells.Item(80, 12).Value2 = 3771.2942  # L80
$ws.Cells.Item(80, 13).Value2 = -1786.074  # M80
$ws.Cells.Item(80, 14).Value2 = -5767.2942  # N80
$ws.Cells.Item(83, 8).Value2 = 3165.5  # H83
$ws.Cells.Item(83, 9).Value2 = 2784.074  # I83
$ws.Cells.Item(83, 10).Value2 = 3771.2942  # J83
$ws.Cells.Item(83, 11).Value2 = 13920.37  # K83
$ws.Cells.Item(83, 12).Value2 = 18856.471  # L83
$ws.Cells.Item(83, 13).Value2 = -8928.370000000001  # M83
$ws.Cells.Item(83, 14).Value2 = -28840.471  # N83
$ws.Cells.Item(102, 8).Value2 = 37503.207  # H102
$ws.Cells.Item(102, 9).Value2 = 2367.158  # I102
$ws.Cells.Item(102, 10).Value2 = 104261.7  # J102
$ws.Cells.Item(102, 11).Value2 = 2367.158  # K102
$ws.Cells.Item(102, 12).Value2 = 104261.7  # L102
$ws.Cells.Item(102, 13).Value2 = -745.1579999999999  # M102
$ws.Cells.Item(102, 14).Value2 = -107505.7  # N102
$ws.Cells.Item(132, 8).Value2 = 2247.8816  # H132
$ws.Cells.Item(132, 9).Value2 = 1851.0862  # I132
$ws.Cells.Item(132, 10).Value2 = 3526.4443  # J132
$ws.Cells.Item(132, 11).Value2 = 5553.2586  # K132
$ws.Cells.Item(132, 12).Value2 = 10579.3329  # L132
$ws.Cells.Item(132, 13).Value2 = -3023.2586  # M132
$ws.Cells.Item(132, 14).Value2 = -15639.3329  # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value2 = 47004  # H4
$ws.Cells.Item(4, 10).Value2 = 47004  # J4
$ws.Cells.Item(4, 12).Value2 = 47004  # L4
$ws.Cells.Item(4, 14).Value2 = -47230  # N4
$ws.Cells.Item(7, 8).Value2 = 2372.9375  # H7
$ws.Cells.Item(7, 9).Value2 = 1843.6154  # I7
$ws.Cells.Item(7, 10).Value2 = 4666.6665  # J7
$ws.Cells.Item(7, 11).Value2 = 1843.6154  # K7
$ws.Cells.Item(7, 12).Value2 = 4666.6665  # L7
$ws.Cells.Item(7, 13).Value2 = -1731.6154  # M7
$ws.Cells.Item(7, 14).Value2 = -4890.6665  # N7
$ws.Cells.Item(28, 8).Value2 = 47004  # H28
$ws.Cells.Item(28, 10).Value2 = 47004  # J28
$ws.Cells.Item(28, 12).Value2 = 47004  # L28
$ws.Cells.Item(28, 14).Value2 = -47468  # N28
$ws.Cells.Item(37, 8).Value2 = 47004  # H37
$ws.Cells.Item(37, 10).Value2 = 47004  # J37
$ws.Cells.Item(37, 12).Value2 = 47004  # L37
$ws.Cells.Item(37, 14).Value2 = -47218  # N37
$ws.Cells.Item(42, 8).Value2 = 80021  # H42
$ws.Cells.Item(42, 10).Value2 = 80021  # J42
$ws.Cells.Item(42, 12).Value2 = 80021  # L42
$ws.Cells.Item(42, 14).Value2 = -81147  # N42
$ws.Cells.Item(49, 8).Value2 = 80021  # H49
$ws.Cells.Item(49, 10).Value2 = 80021  # J49
$ws.Cells.Item(49, 12).Value2 = 80021  # L49
$ws.Cells.Item(49, 14).Value2 = -80315  # N49
$ws.Cells.Item(122, 8).Value2 = 2952.7632  # H122
$ws.Cells.Item(122, 9).Value2 = 2456.8333  # I122
$ws.Cells.Item(122, 11).Value2 = 7370.499899999999  # K122
$ws.Cells.Item(122, 13).Value2 = -4920.499899999999  # M122
$ws.Cells.Item(126, 8).Value2 = 2372.9375  # H126
$ws.Cells.Item(126, 9).Value2 = 1843.6154  # I126
$ws.Cells.Item(126, 10).Value2 = 4666.6665  # J126
$ws.Cells.Item(126, 11).Value2 = 5530.8462  # K126
$ws.Cells.Item(126, 12).Value2 = 13999.9995  # L126
$ws.Cells.Item(126, 13).Value2 = -3060.8462  # M126
$ws.Cells.Item(126, 14).Value2 = -18939.9995  # N126
$ws.Cells.Item(132, 8).Value2 = 2045.695  # H132
$ws.Cells.Item(132, 9).Value2 = 1374  # I132
$ws.Cells.Item(132, 11).Value2 = 4122  # K132
$ws.Cells.Item(132, 13).Value2 = -1592  # M132

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value2 = 0  # H51
$ws.Cells.Item(51, 10).Value2 = 0  # J51
$ws.Cells.Item(51, 12).Value2 = 0  # L51
$ws.Cells.Item(51, 14).ClearContents()  # N51
$ws.Cells.Item(122, 8).Value2 = 1972.5588  # H122
$ws.Cells.Item(122, 9).Value2 = 1598.8966  # I122
$ws.Cells.Item(122, 11).Value2 = 4796.6898  # K122
$ws.Cells.Item(122, 13).Value2 = -2346.6898  # M122
